$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "Julianaju"
$ws.Range("E2").Value = "Julianaju"
$ws.Range("E2").Select()
